# Applies the 61.xlsx (Psalm 61) update:
#  - Adds Hebrew verse-number letters (א, ב, ג) to the first three verse
#    header rows (2, 7, 13) which previously only held the Arabic numeral.
#  - Re-styles every verse-header row (2,7,13,24,32,42,52,61,70) so that
#    column A (the Hebrew letter) is right-aligned Calibri 10, and column B
#    (the Arabic numeral) is right-aligned Calibri 8.
#  - Fixes the Russian translation of "Начальнику хора" -> "Руководителю хора".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that mark the start of a new verse, paired with the Hebrew letter
# that should appear in column A.
$verseRows = @{
    2  = "א"
    7  = "ב"
    13 = "ג"
    24 = "ד"
    32 = "ה"
    42 = "ו"
    52 = "ז"
    61 = "ח"
    70 = "ט"
}

foreach ($row in $verseRows.Keys) {
    $letter = $verseRows[$row]

    $cellA = $ws.Cells.Item($row, 1)
    $cellB = $ws.Cells.Item($row, 2)

    # Rows 2, 7 and 13 only had the bare Arabic numeral in column A before;
    # give them the Hebrew verse letter like the later verses already have.
    if ($row -eq 2 -or $row -eq 7 -or $row -eq 13) {
        $cellA.Value = $letter
    }

    # Column A (Hebrew verse letter): Calibri 10, black, right-aligned.
    $cellA.HorizontalAlignment = -4152
    $cellA.Font.Color = 0
    $cellA.Font.Name = "Calibri"
    $cellA.Font.Size = 10

    # Column B (Arabic verse number): Calibri 8, black, right-aligned.
    $cellB.HorizontalAlignment = -4152
    $cellB.Font.Color = 0
    $cellB.Font.Name = "Calibri"
    $cellB.Font.Size = 8
}

# Correct the translation of "Начальнику хора" -> "Руководителю хора".
$ws.Range("B3").Value = "Руководителю хора"

Write-Host "done"
